$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Rows 62-69: the trip duration ("Duracion_Trayecto_Min", col O) and the
# frequency ("Frecuencia_Min", col P) values had been entered in the wrong
# columns. Swap them back: P's value moves into O, and P is reset to 0.
for ($r = 62; $r -le 69; $r++) {
    $durCell  = $ws.Range("O$r")
    $freqCell = $ws.Range("P$r")

    $oVal = $durCell.Value2
    $pVal = $freqCell.Value2

    $durCell.Value = $pVal
    $freqCell.Value = $oVal
}

# Restore the window/selection state: scrolled down so row 46 is at the top,
# with P63 as the active cell.
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P63").Select()
